$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "296.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.75%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.77"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.029"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.33%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07558"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.57%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.391"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.50%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.591"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.88%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9290"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.82%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1204"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.80%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1839"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.58%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09022"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.25%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03993"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.53%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1052"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.09%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001282"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.64%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005798"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.26%"
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.003961"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "4.57%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.353"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.36%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3320"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.32%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.893"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.67%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1420"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.83%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3000"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.99%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04065"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5.48%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001267"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.03%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-3.92%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.04%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02414"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.70%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05221"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.72%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006035"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-7.69%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007790"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.56%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007539"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.32%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007846"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "10.67%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006792"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.90%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.02%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04643"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "126.83%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004203"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.04%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.02%"
